$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New competition row: "Dog Breeds" (lesson 1, dog breed identification competition)
# Copy formatting from row 2 (the existing competition row) down to row 3 first,
# so the new row picks up the same per-column styles (bold/centered/date/percent).
$ws.Range("A2:L2").Copy() | Out-Null
$ws.Range("A3:L3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A3").Value2 = 2
$ws.Range("L3").Value2 = "https://www.kaggle.com/c/dog-breed-identification"
$ws.Range("B3").Value2 = "Dog Breeds"
$ws.Range("C3").Formula = "=I3/J3"
$ws.Range("D3").Value2 = "Ended"
$ws.Range("E3").Value2 = "Late"
$ws.Range("F3").Value2 = 43201
$ws.Range("G3").Value2 = "Lesson 1"
$ws.Range("H3").Value2 = 0.21451
$ws.Range("I3").Value2 = 311
$ws.Range("J3").Value2 = 1286

# Note column (K) has no entry for this competition - clear the copied format's cell
$ws.Range("K3").Clear() | Out-Null

# Move the active selection down to below the newly added row, as in the source file
$ws.Range("B4").Select() | Out-Null

$wb.Save()
